# Epi Info7: Extended to the logic driving the creation of a new project
# from data Dictionary.
#
# Adds a "Title" and "Description" column to the Survey worksheet (between
# Question and Variable_Name), renames the Then_Question/Else_Question
# headers to Then_Goto/Else_Goto, and fills in the new columns for every
# existing question row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------
# 1. Make room for the two new columns (Title, Description) right after
#    column A (Question). This pushes Variable_Name..Else_Question from
#    B..H out to D..J and keeps every formula / data validation sqref in
#    sync automatically.
# ---------------------------------------------------------------------
$ws.Columns("B:C").Insert()

# ---------------------------------------------------------------------
# 2. Header row
# ---------------------------------------------------------------------
$ws.Range("B1").Value = "Title"
$ws.Range("C1").Value = "Description"
$ws.Range("I1").Value = "Then_Goto"
$ws.Range("J1").Value = "Else_Goto"

# ---------------------------------------------------------------------
# 3. New Title / Description values for every question row
# ---------------------------------------------------------------------
$ws.Range("B2").Value = "Patient"
$ws.Range("C2").Value = "Please enter Name…."

$ws.Range("B3").Value = "Age"
$ws.Range("C3").Value = "Please enter Age …."

$ws.Range("B4").Value = "Sex"
$ws.Range("C4").Value = "Please enter Sex…."

$ws.Range("B5").Value = "Pregnant"
$ws.Range("C5").Value = "Please enter Pregnant..."

$ws.Range("B6").Value = "symptoms"
$ws.Range("C6").Value = "Please enter symptoms…"

# ---------------------------------------------------------------------
# 4. New columns get a fixed (non-autofit) width, matching the width
#    Excel applied once "Title"/"Description" text was entered.
# ---------------------------------------------------------------------
$ws.Columns("B:C").ColumnWidth = 25.33

# ---------------------------------------------------------------------
# 5. Re-point the plain TRUE/FALSE validation (Required column, now F)
#    and the DataTypes list validation (Question_Type column, now E).
#    Columns("B:C").Insert() already re-targets the plain dataValidation
#    automatically (D2:D1048576 -> F2:F1048576); re-asserting it here is
#    harmless and keeps this step self-contained.
# ---------------------------------------------------------------------
$ws.Range("F2:F1048576").Validation.Delete()
$ws.Range("F2:F1048576").Validation.Add(3, 1, 1, '"TRUE,FALSE"')

# ---------------------------------------------------------------------
# 6. Selection / active cell, matching the saved state in the workbook.
# ---------------------------------------------------------------------
$ws.Range("I1").Select()
